# Apply the graphical updates:
#  - Shared-string table for the Park factor levels gets re-sorted
#    alphabetically, and "ParkEiffel" is corrected to "ParkEifel".
#  - Because the row order (3,4,5,8) in both model-summary tables is tied
#    to the old factor-level order, the Value/Std.Error/p numbers in those
#    rows have to move together with their (corrected) label so the
#    correct number still sits next to the correct park name.
#
# Concretely, for rows 3, 4, 5 and 8 on each sheet:
#   new row3 (ParkEifel)        <- old row8 (was ParkEiffel) numbers
#   new row4 (ParkHainich)      <- old row3 numbers
#   new row5 (ParkHunsrueck)    <- old row4 numbers
#   new row8 (ParkSaechs_Schw)  <- old row5 numbers
# Rows 6 (ParkJasmund) and 7 (ParkKellerwald) are unaffected - their
# alphabetical position does not change.

$wb = $excel.ActiveWorkbook

foreach ($ws in @($wb.Worksheets.Item(1), $wb.Worksheets.Item(2))) {

    # Capture the current (pre-edit) numeric columns for the four rows
    # that need to be re-shuffled.
    $b3 = $ws.Range("B3").Value2
    $c3 = $ws.Range("C3").Value2
    $d3 = $ws.Range("D3").Value2

    $b4 = $ws.Range("B4").Value2
    $c4 = $ws.Range("C4").Value2
    $d4 = $ws.Range("D4").Value2

    $b5 = $ws.Range("B5").Value2
    $c5 = $ws.Range("C5").Value2
    $d5 = $ws.Range("D5").Value2

    $b8 = $ws.Range("B8").Value2
    $c8 = $ws.Range("C8").Value2
    $d8 = $ws.Range("D8").Value2

    # Row 3 becomes "ParkEifel" (renamed from "ParkEiffel") and takes the
    # numbers that used to belong to row 8.
    $ws.Range("A3").Value2 = "ParkEifel"
    $ws.Range("B3").Value2 = $b8
    $ws.Range("C3").Value2 = $c8
    $ws.Range("D3").Value2 = $d8

    # Row 4 becomes "ParkHainich" (previously row 3's label) and takes the
    # numbers that used to belong to row 3.
    $ws.Range("A4").Value2 = "ParkHainich"
    $ws.Range("B4").Value2 = $b3
    $ws.Range("C4").Value2 = $c3
    $ws.Range("D4").Value2 = $d3

    # Row 5 becomes "ParkHunsrueck" (previously row 4's label) and takes
    # the numbers that used to belong to row 4.
    $ws.Range("A5").Value2 = "ParkHunsrueck"
    $ws.Range("B5").Value2 = $b4
    $ws.Range("C5").Value2 = $c4
    $ws.Range("D5").Value2 = $d4

    # Row 8 becomes "ParkSaechs_Schw" (previously row 5's label) and takes
    # the numbers that used to belong to row 5.
    $ws.Range("A8").Value2 = "ParkSaechs_Schw"
    $ws.Range("B8").Value2 = $b5
    $ws.Range("C8").Value2 = $c5
    $ws.Range("D8").Value2 = $d5

    # Rows 6 (ParkJasmund) and 7 (ParkKellerwald) keep their label and
    # numbers exactly as they were - nothing to do.
}
